$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1026.037
$ws.Range("I15").Value = 1026.037
$ws.Range("K15").Value = 3078.111
$ws.Range("M15").Value = -2909.111
$ws.Range("H33").Value = 103.31579
$ws.Range("I33").Value = 104.4
$ws.Range("K33").Value = 104.4
$ws.Range("M33").Value = 124.6
$ws.Range("H51").Value = 4727.8823
$ws.Range("I51").Value = 4495
$ws.Range("K51").Value = 4495
$ws.Range("M51").Value = -4011
$ws.Range("H98").Value = 1762.4
$ws.Range("I98").Value = 1506.0667
$ws.Range("K98").Value = 1506.0667
$ws.Range("M98").Value = -8.066700000000083
$ws.Range("H122").Value = 1762.4
$ws.Range("I122").Value = 1506.0667
$ws.Range("K122").Value = 4518.2001
$ws.Range("M122").Value = -2068.2001
$ws.Range("H138").Value = 5349.051
$ws.Range("I138").Value = 1158.5
$ws.Range("J138").Value = 6204.265
$ws.Range("K138").Value = 3475.5
$ws.Range("L138").Value = 18612.795
$ws.Range("M138").Value = 1664.5
$ws.Range("N138").Value = -28892.795
$ws.Range("H141").Value = 10070.714
$ws.Range("I141").Value = 12199
$ws.Range("K141").Value = 36597
$ws.Range("M141").Value = -31417

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3451.1448
$ws.Range("I32").Value = 2157.5688
$ws.Range("K32").Value = 2157.5688
$ws.Range("M32").Value = -1870.5688
$ws.Range("H45").Value = 3576.2
$ws.Range("I45").Value = 3345.25
$ws.Range("K45").Value = 3345.25
$ws.Range("M45").Value = -2968.25
$ws.Range("H74").Value = 13893929
$ws.Range("I74").Value = 35716030
$ws.Range("K74").Value = 35716030
$ws.Range("M74").Value = -35715156
$ws.Range("H77").Value = 13893929
$ws.Range("I77").Value = 35716030
$ws.Range("K77").Value = 178580150
$ws.Range("M77").Value = -178575782
$ws.Range("H132").Value = 29197.885
$ws.Range("I132").Value = 27254.13
$ws.Range("J132").Value = 44100
$ws.Range("K132").Value = 81762.39
$ws.Range("L132").Value = 132300
$ws.Range("M132").Value = -79232.39
$ws.Range("N132").Value = -137360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2326.7693
$ws.Range("I86").Value = 1841.3334
$ws.Range("K86").Value = 1841.3334
$ws.Range("M86").Value = -718.3334
$ws.Range("H89").Value = 2326.7693
$ws.Range("I89").Value = 1841.3334
$ws.Range("K89").Value = 9206.666999999999
$ws.Range("M89").Value = -3590.666999999999
$ws.Range("H107").Value = 15875648
$ws.Range("I107").Value = 20204844
$ws.Range("K107").Value = 20204844
$ws.Range("M107").Value = -20202924
$ws.Range("H134").Value = 4659.727
$ws.Range("I134").Value = 2383.4546
$ws.Range("K134").Value = 7150.3638
$ws.Range("M134").Value = -4615.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17547364
$ws.Range("I31").Value = 35715856
$ws.Range("J31").Value = 5372.793
$ws.Range("K31").Value = 35715856
$ws.Range("L31").Value = 5372.793
$ws.Range("M31").Value = -35715561
$ws.Range("N31").Value = -5962.793
$ws.Range("H34").Value = 17547364
$ws.Range("I34").Value = 35715856
$ws.Range("J34").Value = 5372.793
$ws.Range("K34").Value = 35715856
$ws.Range("L34").Value = 5372.793
$ws.Range("M34").Value = -35715654
$ws.Range("N34").Value = -5776.793
$ws.Range("H58").Value = 1540766.9
$ws.Range("I58").Value = 1668930.8
$ws.Range("J58").Value = 2800
$ws.Range("K58").Value = 1668930.8
$ws.Range("L58").Value = 2800
$ws.Range("M58").Value = -1668727.8
$ws.Range("N58").Value = -3206
$ws.Range("H136").Value = 1540766.9
$ws.Range("I136").Value = 1668930.8
$ws.Range("J136").Value = 2800
$ws.Range("K136").Value = 5006792.4
$ws.Range("L136").Value = 8400
$ws.Range("M136").Value = -5004242.4
$ws.Range("N136").Value = -13500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4625433.5
$ws.Range("J4").Value = 14026000
$ws.Range("L4").Value = 42078000
$ws.Range("N4").Value = -42078224
$ws.Range("H68").Value = 462956.12
$ws.Range("J68").Value = 561732.5
$ws.Range("L68").Value = 1685197.5
$ws.Range("N68").Value = -1686819.5
$ws.Range("H71").Value = 462956.12
$ws.Range("J71").Value = 561732.5
$ws.Range("L71").Value = 5055592.5
$ws.Range("N71").Value = -5063704.5
$ws.Range("H107").Value = 663458.3
$ws.Range("I107").Value = 1154.3334
$ws.Range("J107").Value = 1060840.8
$ws.Range("K107").Value = 3463.0002
$ws.Range("L107").Value = 3182522.4
$ws.Range("M107").Value = -1543.0002
$ws.Range("N107").Value = -3186362.4
$ws.Range("H113").Value = 425.5625
$ws.Range("I113").Value = 283.8889
$ws.Range("K113").Value = 851.6667
$ws.Range("M113").Value = 1318.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1671895.2
$ws.Range("I80").Value = 2776992.2
$ws.Range("J80").Value = 14249.75
$ws.Range("K80").Value = 2776992.2
$ws.Range("L80").Value = 14249.75
$ws.Range("M80").Value = -2775994.2
$ws.Range("N80").Value = -16245.75
$ws.Range("H83").Value = 1671895.2
$ws.Range("I83").Value = 2776992.2
$ws.Range("J83").Value = 14249.75
$ws.Range("K83").Value = 13884961
$ws.Range("L83").Value = 71248.75
$ws.Range("M83").Value = -13879969
$ws.Range("N83").Value = -81232.75
$ws.Range("H107").Value = 2506742
$ws.Range("I107").Value = 3968573.2
$ws.Range("K107").Value = 3968573.2
$ws.Range("M107").Value = -3966653.2
$ws.Range("H122").Value = 309731.1
$ws.Range("I122").Value = 481319.12
$ws.Range("K122").Value = 1443957.36
$ws.Range("M122").Value = -1441507.36
$ws.Range("H132").Value = 3733.1667
$ws.Range("I132").Value = 3699.1765
$ws.Range("J132").Value = 3877.625
$ws.Range("K132").Value = 11097.5295
$ws.Range("L132").Value = 11632.875
$ws.Range("M132").Value = -8567.529500000001
$ws.Range("N132").Value = -16692.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4285.625
$ws.Range("I7").Value = 3776.795
$ws.Range("J7").Value = 5452.9414
$ws.Range("K7").Value = 3776.795
$ws.Range("L7").Value = 5452.9414
$ws.Range("M7").Value = -3664.795
$ws.Range("N7").Value = -5676.9414
$ws.Range("H55").Value = 265.95456
$ws.Range("J55").Value = 93
$ws.Range("L55").Value = 93
$ws.Range("N55").Value = -439
$ws.Range("H100").Value = 2681.7778
$ws.Range("I100").Value = 2361.6667
$ws.Range("J100").Value = 3322
$ws.Range("K100").Value = 2361.6667
$ws.Range("L100").Value = 3322
$ws.Range("M100").Value = -1820.6667
$ws.Range("N100").Value = -4404
$ws.Range("H126").Value = 4285.625
$ws.Range("I126").Value = 3776.795
$ws.Range("J126").Value = 5452.9414
$ws.Range("K126").Value = 11330.385
$ws.Range("L126").Value = 16358.8242
$ws.Range("M126").Value = -8860.385
$ws.Range("N126").Value = -21298.8242
$ws.Range("H136").Value = 3897.36
$ws.Range("I136").Value = 3957.6875
$ws.Range("J136").Value = 2449.5
$ws.Range("K136").Value = 11873.0625
$ws.Range("L136").Value = 7348.5
$ws.Range("M136").Value = -9323.0625
$ws.Range("N136").Value = -12448.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 86409.086
$ws.Range("I107").Value = 114152
$ws.Range("K107").Value = 342456
$ws.Range("M107").Value = -340536
$ws.Range("H126").Value = 3791
$ws.Range("I126").Value = 3250.2
$ws.Range("K126").Value = 9750.599999999999
$ws.Range("M126").Value = -7280.599999999999
$ws.Range("H132").Value = 15436937
$ws.Range("I132").Value = 2651164.8
$ws.Range("J132").Value = 33337018
$ws.Range("K132").Value = 7953494.399999999
$ws.Range("L132").Value = 100011054
$ws.Range("M132").Value = -7950964.399999999
$ws.Range("N132").Value = -100016114
$ws.Range("H136").Value = 6499.18
$ws.Range("I136").Value = 3012.634
$ws.Range("J136").Value = 8922.034
$ws.Range("K136").Value = 9037.902
$ws.Range("L136").Value = 26766.102
$ws.Range("M136").Value = -6487.902
$ws.Range("N136").Value = -31866.102
